$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.123.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "'2.420.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'488.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "'155.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.609"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +18.83%  "
$ws.Range("D9").Value = "'2.448.74"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("E10").Value = "  +10.02%  "
$ws.Range("D11").Value = "'0.100"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").Value = "'2.843.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").Value = "'57.198.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "'20.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("E17").Value = "  -3.24%  "
$ws.Range("D18").Value = "'2.451.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("D20").Value = "'323.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("D21").Value = "'10.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").Value = "'58.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "'0.402"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").Value = "'2.546.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").Value = "'7.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.57%  "
$ws.Range("D30").Value = "'0.0₃0796"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.91%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "'151.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "'18.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.00%  "
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("E35").Value = "  +1.80%  "
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").Value = "'0.827"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.41%  "
$ws.Range("D39").Value = "'34.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").Value = "'3.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("E42").Value = "  +4.41%  "
$ws.Range("D43").Value = "'0.996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "'276.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.47%  "
$ws.Range("E46").Value = "  -4.43%  "
$ws.Range("D47").Value = "'10.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "'0.0229"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("D49").Value = "'4.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.78%  "
$ws.Range("D50").Value = "'17.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").Value = "'0.679"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.07%  "
